$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 3574
$ws.Range("F4").Value = 379
$ws.Range("F5").Value = 8341
$ws.Range("F7").Value = 130
$ws.Range("F8").Value = 2240
$ws.Range("F10").Value = 105
$ws.Range("F11").Value = 76
$ws.Range("F12").Value = 660
$ws.Range("F13").Value = 112
$ws.Range("F14").Value = 7503
$ws.Range("F16").Value = 7725
$ws.Range("F18").Value = 57893
$ws.Range("F19").Value = 57893
$ws.Range("F20").Value = 4828
$ws.Range("F21").Value = 1064
$ws.Range("F22").Value = 957
$ws.Range("F23").Value = 511
$ws.Range("F25").Value = 932
$ws.Range("F28").Value = 5315
$ws.Range("F30").Value = 122
$ws.Range("F33").Value = 1411
$ws.Range("F34").Value = 1990
$ws.Range("F40").Value = 732
$ws.Range("F41").Value = 43
$ws.Range("F43").Value = 275
$ws.Range("F44").Value = 248
$ws.Range("F48").Value = 16
$ws.Range("F50").Value = 2490

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 204
$ws.Range("F6").Value = 155
$ws.Range("F9").Value = 7660
$ws.Range("F12").Value = 6
$ws.Range("F14").Value = 7
$ws.Range("F20").Value = 23
$ws.Range("F24").Value = 41
$ws.Range("F41").Value = 7

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 2382
$ws.Range("F5").Value = 1619
$ws.Range("F9").Value = 9448
$ws.Range("F10").Value = 1773
$ws.Range("F11").Value = 182
$ws.Range("F12").Value = 115
$ws.Range("F16").Value = 2417
$ws.Range("F17").Value = 140
$ws.Range("F18").Value = 68
$ws.Range("F19").Value = 532

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 3575
$ws.Range("F4").Value = 2382
$ws.Range("F5").Value = 8342
$ws.Range("F6").Value = 1773
$ws.Range("F7").Value = 182
$ws.Range("F9").Value = 140
$ws.Range("F10").Value = 76
$ws.Range("F11").Value = 660
$ws.Range("F12").Value = 7725
$ws.Range("F13").Value = 57893
$ws.Range("F14").Value = 204
$ws.Range("F16").Value = 4828
$ws.Range("F17").Value = 957
$ws.Range("F18").Value = 511
$ws.Range("F19").Value = 932
$ws.Range("F21").Value = 155
$ws.Range("F22").Value = 122
$ws.Range("F24").Value = 1411
$ws.Range("F25").Value = 1990
$ws.Range("F27").Value = 532
$ws.Range("F29").Value = 7
$ws.Range("F33").Value = 23
$ws.Range("F34").Value = 43
$ws.Range("F36").Value = 275
$ws.Range("F47").Value = 2490
